# Update the "Last Modified" style auto date placeholders (on every slide
# layout, the slide master, and the notes master) from 11/21/2019 to the
# new cached value 11/22/19, and recenter/resize the dataset screenshot
# picture on slide 5.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Type -eq 14) {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $shape.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$newDate = "11/22/19"

# Slide master.
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Notes master.
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# Slide 5: recenter/resize the dataset-source picture.
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shape = $slide5.Shapes.Item($i)
    if ($shape.Name -eq "Picture 4") {
        $shape.Left = 96.778702
        $shape.Top = 211.281692
        $shape.Width = 797.570038
        $shape.Height = 183.498932
    }
}
